# Mise à jour du planning
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Corriger le libellé de la tâche "Déterminer la liste des étapes à suivre"
# (retrait de la mention "(pseudo code)")
$ws.Range("B11").Value = "Déterminer la liste des étapes à suivre  pour arriver à la solution"

# Les deux premières sous-tâches de "Contexte et état de l'art" sont terminées :
# passage du statut "En cours" (jaune) à "Terminé" (vert), comme B7.
$ws.Range("B6").Interior.Color = 5287936
$ws.Range("B8").Interior.Color = 5287936

# Mettre à jour la sélection affichée et remonter la vue en haut de la feuille
$ws.Range("A6:B8").Select()
$excel.ActiveWindow.ScrollRow = 1
